$d = $word.ActiveDocument

$replacements = @(
    @{old="437×2=874"; new="887×8=7096"},
    @{old="639×3=1917"; new="953×5=4765"},
    @{old="729×8=5832"; new="463×4=1852"},
    @{old="626×8=5008"; new="141×4=564"},
    @{old="222×2=444"; new="579×5=2895"},
    @{old="548×2=1096"; new="169×5=845"},
    @{old="597×6=3582"; new="965×7=6755"},
    @{old="788×5=3940"; new="407×4=1628"},
    @{old="342×9=3078"; new="195×8=1560"},
    @{old="772×3=2316"; new="186×4=744"},
    @{old="734×6=4404"; new="807×7=5649"},
    @{old="225×3=675"; new="479×9=4311"},
    @{old="741×4=2964"; new="566×8=4528"},
    @{old="621×4=2484"; new="857×4=3428"},
    @{old="383×5=1915"; new="716×8=5728"},
    @{old="967×6=5802"; new="794×5=3970"},
    @{old="306×9=2754"; new="906×8=7248"},
    @{old="554×2=1108"; new="887×5=4435"},
    @{old="261×6=1566"; new="616×2=1232"},
    @{old="984×7=6888"; new="868×9=7812"},
    @{old="904×7=6328"; new="944×5=4720"},
    @{old="445×6=2670"; new="632×9=5688"},
    @{old="511×3=1533"; new="133×7=931"},
    @{old="594×9=5346"; new="859×2=1718"},
    @{old="473×2=946"; new="407×5=2035"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
